$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells ---
$ws.Range("B1").Value = "2025-07-10 11:36:59"
$ws.Range("B2").Value = "PUne"

# --- Row 5 (existing row, updated values) ---
$ws.Range("A5").Value = " Afsar"
$ws.Range("B5").Value = " REDX"
$ws.Range("C5").Value = " Pune"
$ws.Range("D5").Value = " Mr.Sachin"
$ws.Range("E5").Value = "2025-07-10 05:40:09"
$ws.Range("F5").Value = "Honeywell"
$ws.Range("G5").Value = "5MP (2560x1920)"
$ws.Range("H5").Value = 25
$ws.Range("I5").Value = "H265"
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 30
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.17
$ws.Range("N5").Value = 0.06

# --- Row 6 (existing row, updated values) ---
$ws.Range("A6").Value = "Afsar Khan"
$ws.Range("B6").Value = "REDX_VMS Portal"
$ws.Range("C6").Value = "Pune, Maharashtra India"
$ws.Range("D6").Value = "Mr.Sahil"
$ws.Range("E6").Value = "2025-06-30 13:24:28"
$ws.Range("F6").Value = "Honeywell"
$ws.Range("G6").Value = "12MP (4000x3000)"
$ws.Range("H6").Value = 25
$ws.Range("I6").Value = "H265"
$ws.Range("J6").Value = 24
$ws.Range("K6").Value = 30
$ws.Range("L6").Value = 99
$ws.Range("M6").Value = 1421.35
$ws.Range("N6").Value = 460.52

# --- Row 7 (new row) ---
$ws.Range("A7").Value = "Afsar Khan"
$ws.Range("B7").Value = "REDX_VMS Portal"
$ws.Range("C7").Value = "Pune, Maharashtra India"
$ws.Range("D7").Value = "Mr.Sahil"
$ws.Range("E7").Value = "2025-06-30 13:24:28"
$ws.Range("F7").Value = "Honeywell"
$ws.Range("G7").Value = "5MP (2560x1920)"
$ws.Range("H7").Value = 25
$ws.Range("I7").Value = "H265"
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 30
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1421.35
$ws.Range("N7").Value = 460.52

# --- Row 8 (new row) ---
$ws.Range("A8").Value = "Afsar Khan"
$ws.Range("B8").Value = "REDX_VMS Portal"
$ws.Range("C8").Value = "Pune, Maharashtra India"
$ws.Range("D8").Value = "Mr.Sahil"
$ws.Range("E8").Value = "2025-06-30 13:24:28"
$ws.Range("F8").Value = "Honeywell"
$ws.Range("G8").Value = "5MP (2560x1920)"
$ws.Range("H8").Value = 25
$ws.Range("I8").Value = "H265"
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 30
$ws.Range("L8").Value = 42
$ws.Range("M8").Value = 1421.35
$ws.Range("N8").Value = 460.52
